# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 12-13 (pushing the existing
# rows 12-22 down to rows 14-24), then fill the two new rows with the
# new data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the current row 12; this shifts existing rows
# 12:22 down to 14:24, carrying their contents (and column D's date
# style) with them.
$ws.Rows("12:13").Insert()

# New row 12
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44579
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 100112030
$ws.Cells.Item(12, 7).Value = "Poroto granado"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 50
$ws.Cells.Item(12, 11).Value = 35000
$ws.Cells.Item(12, 12).Value = 35000
$ws.Cells.Item(12, 13).Value = 35000
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 1400
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# New row 13
$ws.Cells.Item(13, 1).Value = 4
$ws.Cells.Item(13, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(13, 3).Value = "Los Lagos"
$ws.Cells.Item(13, 4).Value = 44579
$ws.Cells.Item(13, 5).Value = 10
$ws.Cells.Item(13, 6).Value = 100112030
$ws.Cells.Item(13, 7).Value = "Poroto granado"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 30000
$ws.Cells.Item(13, 12).Value = 30000
$ws.Cells.Item(13, 13).Value = 30000
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 1200
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
